$d = $word.ActiveDocument

$d.Content.Find.Execute("41×15=615", $true, $false, $false, $false, $false, $true, 1, $false, "29×16=464", 2) | Out-Null
$d.Content.Find.Execute("23×14=322", $true, $false, $false, $false, $false, $true, 1, $false, "70×76=5320", 2) | Out-Null
$d.Content.Find.Execute("48×63=3024", $true, $false, $false, $false, $false, $true, 1, $false, "61×98=5978", 2) | Out-Null
$d.Content.Find.Execute("91×41=3731", $true, $false, $false, $false, $false, $true, 1, $false, "28×56=1568", 2) | Out-Null
$d.Content.Find.Execute("53×53=2809", $true, $false, $false, $false, $false, $true, 1, $false, "59×77=4543", 2) | Out-Null
$d.Content.Find.Execute("60×64=3840", $true, $false, $false, $false, $false, $true, 1, $false, "14×70=980", 2) | Out-Null
$d.Content.Find.Execute("84×18=1512", $true, $false, $false, $false, $false, $true, 1, $false, "59×99=5841", 2) | Out-Null
$d.Content.Find.Execute("24×29=696", $true, $false, $false, $false, $false, $true, 1, $false, "50×24=1200", 2) | Out-Null
$d.Content.Find.Execute("70×44=3080", $true, $false, $false, $false, $false, $true, 1, $false, "83×90=7470", 2) | Out-Null
$d.Content.Find.Execute("48×55=2640", $true, $false, $false, $false, $false, $true, 1, $false, "55×81=4455", 2) | Out-Null
$d.Content.Find.Execute("16×54=864", $true, $false, $false, $false, $false, $true, 1, $false, "16×35=560", 2) | Out-Null
$d.Content.Find.Execute("89×35=3115", $true, $false, $false, $false, $false, $true, 1, $false, "90×84=7560", 2) | Out-Null
$d.Content.Find.Execute("96×90=8640", $true, $false, $false, $false, $false, $true, 1, $false, "83×78=6474", 2) | Out-Null
$d.Content.Find.Execute("68×92=6256", $true, $false, $false, $false, $false, $true, 1, $false, "42×51=2142", 2) | Out-Null
$d.Content.Find.Execute("41×75=3075", $true, $false, $false, $false, $false, $true, 1, $false, "62×28=1736", 2) | Out-Null
$d.Content.Find.Execute("17×56=952", $true, $false, $false, $false, $false, $true, 1, $false, "76×51=3876", 2) | Out-Null
$d.Content.Find.Execute("67×19=1273", $true, $false, $false, $false, $false, $true, 1, $false, "24×92=2208", 2) | Out-Null
$d.Content.Find.Execute("25×90=2250", $true, $false, $false, $false, $false, $true, 1, $false, "25×77=1925", 2) | Out-Null
$d.Content.Find.Execute("95×98=9310", $true, $false, $false, $false, $false, $true, 1, $false, "78×95=7410", 2) | Out-Null
$d.Content.Find.Execute("60×49=2940", $true, $false, $false, $false, $false, $true, 1, $false, "75×78=5850", 2) | Out-Null
$d.Content.Find.Execute("17×17=289", $true, $false, $false, $false, $false, $true, 1, $false, "54×73=3942", 2) | Out-Null
$d.Content.Find.Execute("54×25=1350", $true, $false, $false, $false, $false, $true, 1, $false, "67×35=2345", 2) | Out-Null
$d.Content.Find.Execute("64×78=4992", $true, $false, $false, $false, $false, $true, 1, $false, "22×38=836", 2) | Out-Null
$d.Content.Find.Execute("34×88=2992", $true, $false, $false, $false, $false, $true, 1, $false, "72×22=1584", 2) | Out-Null
$d.Content.Find.Execute("16×65=1040", $true, $false, $false, $false, $false, $true, 1, $false, "65×76=4940", 2) | Out-Null
